# feat: correct mmcd computation
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) New column H ("fixed mu=0.6" gamma column) for the first block (rows 2-13)
# ---------------------------------------------------------------------------
$ws.Range("H2").Value  = 0.56471816283924903
$ws.Range("H3").Value  = 0.54237288135593198
$ws.Range("H4").Value  = 0.50318471337579596
$ws.Range("H5").Value  = 0.446126447016919
$ws.Range("H6").Value  = 0.54915254237288103
$ws.Range("H7").Value  = 0.45859872611465002
$ws.Range("H8").Value  = 0.33748886910062298
$ws.Range("H9").Value  = 0.42901878914405001
$ws.Range("H10").Value = 0.88535031847133805
$ws.Range("H11").Value = 0.34995547640249303
$ws.Range("H12").Value = 0.41544885177453
$ws.Range("H13").Value = 0.87457627118644099

# relative-performance block (rows 15-26): H = H(row+-13) - $B(row-13)
$ws.Range("H15").Formula = '=H2-$B2'
$ws.Range("H16").Formula = '=H3-$B3'
$ws.Range("H17").Formula = '=H4-$B4'
$ws.Range("H18").Formula = '=H5-$B5'
$ws.Range("H19").Formula = '=H6-$B6'
$ws.Range("H20").Formula = '=H7-$B7'
$ws.Range("H21").Formula = '=H8-$B8'
$ws.Range("H22").Formula = '=H9-$B9'
$ws.Range("H23").Formula = '=H10-$B10'
$ws.Range("H24").Formula = '=H11-$B11'
$ws.Range("H25").Formula = '=H12-$B12'
$ws.Range("H26").Formula = '=H13-$B13'

# Column H inherits the sheet's default column style (percent format, no
# fill) for unstyled rows, but rows 16-20,23,24,26 carry an explicit
# row-level custom format (colored fill) that a brand-new cell would
# otherwise inherit. Force those H cells back to the plain "style 1"
# look (copy format only) to match the rest of the column.
$ws.Range("C1").Copy() | Out-Null
$ws.Range("H16:H20").PasteSpecial(-4122)
$ws.Range("H23").PasteSpecial(-4122)
$ws.Range("H24").PasteSpecial(-4122)
$ws.Range("H26").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# row 28 average: extend the existing per-column average formula into H
$ws.Range("H28").Formula = '=SUM(H15:H26)/12'

# ---------------------------------------------------------------------------
# 2) USOS_MNIST 05/07 block: add a third ("gamma=0.01") comparison column D
# ---------------------------------------------------------------------------
$ws.Range("D60").Value = 0.72222222222222199
$ws.Range("D61").Value = 0.62549999999999994

$ws.Range("D64").Formula = '=D60-$B60'
$ws.Range("D65").Formula = '=D61-$B61'

# ---------------------------------------------------------------------------
# 3) New COIL block at the bottom of the sheet
# ---------------------------------------------------------------------------
$ws.Range("A92").Value = "COIL"
$ws.Range("B92").Value = "MEDA"
$ws.Range("C92").Value = "gamma=0.1"
$ws.Range("D92").Value = "gamma=0.01"

$ws.Range("A93").Value = "COIL_01"
$ws.Range("B93").Value = 0.90138888888888902
$ws.Range("C93").Value = 0.86805555555555602
$ws.Range("D93").Value = 0.88472222222222197
$ws.Range("E93").Value = 0.88888888888888895
$ws.Range("F93").Value = 0.89861111111111103
$ws.Range("G93").Value = 0.86388888888888904

$ws.Range("A94").Value = "COIL_02"
$ws.Range("B94").Value = 0.87083333333333302
$ws.Range("C94").Value = 0.85972222222222205
$ws.Range("D94").Value = 0.87083333333333302
$ws.Range("E94").Value = 0.87083333333333302
$ws.Range("F94").Value = 0.87083333333333302
$ws.Range("G94").Value = 0.85416666666666696

$ws.Range("C96").Formula = '=C93-$B93'
$ws.Range("D96").Formula = '=D93-$B93'
$ws.Range("E96").Formula = '=E93-$B93'
$ws.Range("F96").Formula = '=F93-$B93'
$ws.Range("G96").Formula = '=G93-$B93'

$ws.Range("C97").Formula = '=C94-B94'
$ws.Range("D97").Formula = '=D94-$B94'
$ws.Range("E97").Formula = '=E94-$B94'
$ws.Range("F97").Formula = '=F94-$B94'
$ws.Range("G97").Formula = '=G94-$B94'

$ws.Range("B98").Value = "平均"

$ws.Range("C99").Formula = '=AVERAGE(C96:C97)'
$ws.Range("D99").Formula = '=AVERAGE(D96:D97)'
$ws.Range("E99").Formula = '=AVERAGE(E96:E97)'
$ws.Range("F99").Formula = '=AVERAGE(F96:F97)'
$ws.Range("G99").Formula = '=AVERAGE(G96:G97)'

# ---------------------------------------------------------------------------
# 4) Restore the on-screen selection to what the author left it at
# ---------------------------------------------------------------------------
$ws.Range("F64").Select()
